{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2023-11-11 Saturday\", \"2023-11-12 Sunday\"],\n  [\"32\u00d726=832\", \"36\u00d749=1764\"],\n  [\"22\u00d755=1210\", \"21\u00d730=630\"],\n  [\"51\u00d735=1785\", \"27\u00d760=1620\"],\n  [\"15\u00d752=780\", \"73\u00d799=7227\"],\n  [\"96\u00d731=2976\", \"67\u00d777=5159\"],\n  [\"99\u00d786=8514\", \"50\u00d714=700\"],\n  [\"13\u00d738=494\", \"20\u00d799=1980\"],\n  [\"67\u00d748=3216\", \"28\u00d779=2212\"],\n  [\"50\u00d785=4250\", \"62\u00d785=5270\"],\n  [\"63\u00d721=1323\", \"86\u00d748=4128\"],\n  [\"25\u00d764=1600\", \"49\u00d781=3969\"],\n  [\"83\u00d733=2739\", \"45\u00d761=2745\"],\n  [\"50\u00d751=2550\", \"29\u00d718=522\"],\n  [\"67\u00d720=1340\", \"58\u00d795=5510\"],\n  [\"11\u00d752=572\", \"69\u00d722=1518\"],\n  [\"45\u00d728=1260\", \"55\u00d736=1980\"],\n  [\"26\u00d729=754\", \"53\u00d781=4293\"],\n  [\"39\u00d751=1989\", \"68\u00d736=2448\"],\n  [\"14\u00d747=658\", \"54\u00d742=2268\"],\n  [\"92\u00d760=5520\", \"78\u00d770=5460\"],\n  [\"79\u00d778=6162\", \"49\u00d741=2009\"],\n  [\"81\u00d771=5751\", \"12\u00d765=780\"],\n  [\"31\u00d716=496\", \"14\u00d736=504\"],\n  [\"30\u00d768=2040\", \"64\u00d765=4160\"],\n  [\"87\u00d769=6003\", \"37\u00d763=2331\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-11-11 Saturday\", \"2023-11-12 Sunday\"),\n    @(\"32\u00d726=832\", \"36\u00d749=1764\"),\n    @(\"22\u00d755=1210\", \"21\u00d730=630\"),\n    @(\"51\u00d735=1785\", \"27\u00d760=1620\"),\n    @(\"15\u00d752=780\", \"73\u00d799=7227\"),\n    @(\"96\u00d731=2976\", \"67\u00d777=5159\"),\n    @(\"99\u00d786=8514\", \"50\u00d714=700\"),\n    @(\"13\u00d738=494\", \"20\u00d799=1980\"),\n    @(\"67\u00d748=3216\", \"28\u00d779=2212\"),\n    @(\"50\u00d785=4250\", \"62\u00d785=5270\"),\n    @(\"63\u00d721=1323\", \"86\u00d748=4128\"),\n    @(\"25\u00d764=1600\", \"49\u00d781=3969\"),\n    @(\"83\u00d733=2739\", \"45\u00d761=2745\"),\n    @(\"50\u00d751=2550\", \"29\u00d718=522\"),\n    @(\"67\u00d720=1340\", \"58\u00d795=5510\"),\n    @(\"11\u00d752=572\", \"69\u00d722=1518\"),\n    @(\"45\u00d728=1260\", \"55\u00d736=1980\"),\n    @(\"26\u00d729=754\", \"53\u00d781=4293\"),\n    @(\"39\u00d751=1989\", \"68\u00d736=2448\"),\n    @(\"14\u00d747=658\", \"54\u00d742=2268\"),\n    @(\"92\u00d760=5520\", \"78\u00d770=5460\"),\n    @(\"79\u00d778=6162\", \"49\u00d741=2009\"),\n    @(\"81\u00d771=5751\", \"12\u00d765=780\"),\n    @(\"31\u00d716=496\", \"14\u00d736=504\"),\n    @(\"30\u00d768=2040\", \"64\u00d765=4160\"),\n    @(\"87\u00d769=6003\", \"37\u00d763=2331\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $findText, $true, $false, $false, $false, $false, $true, 1, $false,\n        $replaceText, 2\n    )\n}\n"}
